$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-CellText 'D2' '69.240.18'
Set-CellText 'E2' '  +2.34%  '
Set-CellText 'D3' '3.387.85'
Set-CellText 'E3' '  +1.69%  '
Set-CellText 'E4' '  +0.10%  '
Set-CellText 'D5' '586.86'
Set-CellText 'E5' '  +1.19%  '
Set-CellText 'D6' '180.49'
Set-CellText 'E6' '  +2.82%  '
Set-CellText 'D7' '1.00'
Set-CellText 'E7' '  -0.02%  '
Set-CellText 'E8' '  +1.44%  '
Set-CellText 'E9' '  +7.95%  '
Set-CellText 'D10' '0.593'
Set-CellText 'E10' '  +2.49%  '
Set-CellText 'D11' '48.62'
Set-CellText 'E11' '  +3.91%  '
Set-CellText 'E12' '  +4.08%  '
Set-CellText 'D13' '679.68'
Set-CellText 'E13' '  -1.72%  '
Set-CellText 'D14' '8.65'
Set-CellText 'E14' '  +2.79%  '
Set-CellText 'D15' '3.933.45'
Set-CellText 'E15' '  +1.58%  '
Set-CellText 'D16' '69.330.83'
Set-CellText 'E16' '  +2.44%  '
Set-CellText 'D17' '3.393.68'
Set-CellText 'E17' '  +1.71%  '
Set-CellText 'E18' '  +1.67%  '
Set-CellText 'D19' '17.73'
Set-CellText 'E19' '  +0.71%  '
Set-CellText 'D20' '11.31'
Set-CellText 'E20' '  +2.49%  '
Set-CellText 'E21' '  +1.52%  '
Set-CellText 'E22' '  -0.42%  '
Set-CellText 'D23' '17.16'
Set-CellText 'E23' '  +1.60%  '
Set-CellText 'E24' '  +1.14%  '
Set-CellText 'D25' '3.92'
Set-CellText 'E25' '  +0.47%  '
Set-CellText 'E26' '  +1.92%  '
Set-CellText 'D27' '9.62'
Set-CellText 'E27' '  +1.92%  '
Set-CellText 'D28' '33.89'
Set-CellText 'E28' '  +2.94%  '
Set-CellText 'D29' '8.77'
Set-CellText 'E29' '  +2.73%  '
Set-CellText 'D30' '6.96'
Set-CellText 'E30' '  -1.16%  '
Set-CellText 'D31' '11.15'
Set-CellText 'E31' '  +1.46%  '
Set-CellText 'B32' 'Bittensor'
Set-CellText 'C32' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-CellText 'D32' '555.82'
Set-CellText 'E32' '  -2.28%  '
Set-CellText 'B33' 'dogwifhat'
Set-CellText 'C33' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-CellText 'D33' '3.61'
Set-CellText 'E33' '  +10.30%  '
Set-CellText 'D34' '0.106'
Set-CellText 'E34' '  +1.30%  '
Set-CellText 'D35' '58.57'
Set-CellText 'E35' '  +2.47%  '
Set-CellText 'E36' '  +0.14%  '
Set-CellText 'D37' '3.670.83'
Set-CellText 'E37' '  -1.06%  '
Set-CellText 'E38' '  +4.68%  '
Set-CellText 'D39' '35.57'
Set-CellText 'E39' '  +1.08%  '
Set-CellText 'D40' '0.0₃0718'
Set-CellText 'E40' '  +7.12%  '
Set-CellText 'E41' '  +3.48%  '
Set-CellText 'E42' '  +2.64%  '
Set-CellText 'D43' '0.339'
Set-CellText 'E43' '  +1.46%  '
Set-CellText 'E44' '  +3.87%  '
Set-CellText 'D45' '3.31'
Set-CellText 'E45' '  -0.96%  '
Set-CellText 'D46' '2.68'
Set-CellText 'E46' '  +1.29%  '
Set-CellText 'E47' '  +1.30%  '
Set-CellText 'E48' '  +5.45%  '
Set-CellText 'E49' '  +0.06%  '
Set-CellText 'D50' '133.58'
Set-CellText 'E50' '  +1.12%  '
Set-CellText 'E51' '  +4.77%  '
